$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 74, shifting existing rows 74:107 down to 75:108
$ws.Rows("74:74").Insert()

# Populate the newly inserted row 74 with the new weekly record
$ws.Cells.Item(74, 1).Value = 10
$ws.Cells.Item(74, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(74, 3).Value = "La Araucanía"
$ws.Cells.Item(74, 4).Value = 44529
$ws.Cells.Item(74, 5).Value = 9
$ws.Cells.Item(74, 6).Value = 100112012
$ws.Cells.Item(74, 7).Value = "Espinaca"
$ws.Cells.Item(74, 8).Value = "Sin especificar"
$ws.Cells.Item(74, 9).Value = "Primera"
$ws.Cells.Item(74, 10).Value = 65
$ws.Cells.Item(74, 11).Value = 10000
$ws.Cells.Item(74, 12).Value = 10000
$ws.Cells.Item(74, 13).Value = 10000
$ws.Cells.Item(74, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(74, 15).Value = "Región Metropolitana"
$ws.Cells.Item(74, 16).Value = 1000
$ws.Cells.Item(74, 17).Value = 10
$ws.Cells.Item(74, 18).Value = "Hortaliza"
